# Add a new "timepoint" column (tp / t0 / t5 / t9) to the carbonate-chemistry
# data sheet, inserted right after the "day" column (old column F, now G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F - everything from the old F column onward
# (pH_NBS, DIC_mmolkgSW, TA_umolkgSW, pCO2_matm) shifts one column to the right.
$ws.Columns.Item(6).Insert()

# Header for the new column.
$ws.Range("F1").Value = "tp"

# Rows 2-4 are the t0 (day 0) samples, rows 5-16 are t5 (day 15) samples,
# and rows 17-28 are t9 (day 27) samples.
$ws.Range("F2:F4").Value = "t0"
$ws.Range("F5:F16").Value = "t5"
$ws.Range("F17:F28").Value = "t9"

# Keep column F's width in line with its neighbour (column E).
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Leave the selection where the author's last edit was (new column, last block).
$ws.Range("F17:F28").Select()
